$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Run 50" column (AZ). Excel shifts the remaining columns
# (the old "Mean" column, BA) left into AZ, and the dimension/spans
# shrink from A1:BA14 to A1:AZ14 automatically.
$ws.Columns("AZ").Delete()

# Header: "Gen" -> "MaxFES"
$ws.Range("A1").Value = "MaxFES"

# Column A: generation counts -> fraction-of-MaxFES values
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 0.001
$ws.Range("A4").Value = 0.01
$ws.Range("A5").Value = 0.1
$ws.Range("A6").Value = 0.2
$ws.Range("A7").Value = 0.3
$ws.Range("A8").Value = 0.4
$ws.Range("A9").Value = 0.5
$ws.Range("A10").Value = 0.6
$ws.Range("A11").Value = 0.7
$ws.Range("A12").Value = 0.8
$ws.Range("A13").Value = 0.9
$ws.Range("A14").Value = 1

# Recompute the "Mean" column (now AZ, after the shift above) without
# the removed Run 50 values.
$ws.Range("AZ2").Value = 208.33194889
$ws.Range("AZ3").Value = 171.00440485
$ws.Range("AZ4").Value = 48.91890259
$ws.Range("AZ5").Value = 0.57801372
$ws.Range("AZ6").Value = 0.24442651
$ws.Range("AZ7").Value = 0.16427764
$ws.Range("AZ8").Value = 0.1313659
$ws.Range("AZ9").Value = 0.10636621
$ws.Range("AZ10").Value = 0.08862969
$ws.Range("AZ11").Value = 0.07310155
$ws.Range("AZ12").Value = 0.06635592999999999
$ws.Range("AZ13").Value = 0.05852669
$ws.Range("AZ14").Value = 0.0531551
